# Insert a new record (row 34) into the "Poroto granado" price series.
# All existing rows from 34 downward shift one row lower (34->35, ..., 122->123),
# and row 34 is populated with the newly reported price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(34).Insert()

$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 45260
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 100112030
$ws.Range("G34").Value = "Poroto granado"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 15
$ws.Range("K34").Value = 60000
$ws.Range("L34").Value = 60000
$ws.Range("M34").Value = 60000
$ws.Range("N34").Value = "$/saco 25 kilos"
$ws.Range("O34").Value = "Provincia de Limarí"
$ws.Range("P34").Value = 2400
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"
